$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 76 contains "Subversion Corporation" / "http://subversion.org/" - remove it entirely,
# shifting rows below it up by one (this also shrinks the table "data2" from A1:H90 to A1:H89,
# and drops the two now-unused shared strings from sharedStrings.xml).
$ws.Rows.Item(76).Delete()

# Update the selection to match the post-edit state captured in the diff.
$ws.Range("B74").Select()
